# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Septiembre de 2020 a las 04:34"

# 2. Refresh case numbers for Bolivia (row 30)
$ws.Range("B30").Value = 130470
$ws.Range("C30").Value = 419
$ws.Range("D30").Value = 89032
$ws.Range("E30").Value = 33852
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 36
$ws.Range("H30").Value = 7586

# 3. Refresh case numbers for Australia (row 78)
$ws.Range("B78").Value = 26897
$ws.Range("C78").Value = 12
$ws.Range("D78").Value = 23967
$ws.Range("E78").Value = 2081
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 5
$ws.Range("H78").Value = 849

# 4. Refresh case numbers for Nueva Zelanda (row 156)
$ws.Range("B156").Value = 1815
$ws.Range("C156").Value = 4
$ws.Range("D156").Value = 1719
$ws.Range("E156").Value = 71
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 25

# 5. Timor Oriental / Santa Lucia swap places (rows 204-205)
#    Both rows share identical case figures, only the country names trade rows.
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("B204").Value = 27
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 26
$ws.Range("E204").Value = 1
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0

$ws.Range("A205").Value = "Timor Oriental"
$ws.Range("B205").Value = 27
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 26
$ws.Range("E205").Value = 1
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

# 6. Islas Malvinas / Montserrat swap places (rows 214-215), figures travel with the country
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0
